$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FlowSegment")
$ws.Range("M3").Value = 1935
$ws.Range("O4").Value = 50
$ws.Range("M5").Value = 1440
$ws.Range("M6").Value = 1440
$ws.Range("M7").Value = 1440
$ws.Range("M8").Value = 1440
$ws.Range("M11").Value = 2835
$ws.Range("M12").Value = 2835
$ws.Range("M13").Value = 2835
$ws.Range("M14").Value = 2835
$ws.Range("M15").Value = 2835
$ws.Range("M16").Value = 2835
$ws.Range("M17").Value = 1440
$ws.Range("M18").Value = 1440
$ws.Range("M21").Value = 1440
$ws.Range("M22").Value = 1440
$ws.Range("M23").Value = 1440
$ws.Range("M24").Value = 1440
$ws.Range("M25").Value = 1440
$ws.Range("M26").Value = 1440
$ws.Range("M27").Value = 2835
$ws.Range("M29").Value = 2835
$ws.Range("M30").Value = 2835
$ws.Range("M31").Value = 2835
$ws.Range("M32").Value = 1440
$ws.Range("M33").Value = 2835
$ws.Range("M34").Value = 1440
$ws.Range("O36").Value = 50
$ws.Range("O37").Value = 50
$ws.Range("O38").Value = 50
$ws.Range("O39").Value = 50
$ws.Range("M42").Value = 1440
$ws.Range("M43").Value = 1440
$ws.Range("M45").Value = 2835
$ws.Range("M46").Value = 2835
$ws.Range("M49").Value = 2835
$ws.Range("M50").Value = 2835
$ws.Range("M51").Value = 2835
$ws.Range("M54").Value = 2835
$ws.Range("M56").Value = 2835

$ws = $wb.Worksheets.Item("FlowFitting")
$ws.Range("M52").Value = 1440

$ws = $wb.Worksheets.Item("BuildingElementProxy")
$ws.Range("M4").Value = 1440
$ws.Range("M5").Value = 2835
$ws.Range("M7").Value = 2835
$ws.Range("M11").Value = 1935
$ws.Range("M12").Value = 2835
$ws.Range("M13").Value = 2835
$ws.Range("M22").Value = 2835
$ws.Range("M23").Value = 2835
$ws.Range("M24").Value = 2835
$ws.Range("M25").Value = 2835
$ws.Range("M26").Value = 2835
$ws.Range("M27").Value = 2835
$ws.Range("M28").Value = 2835
$ws.Range("M29").Value = 2835
$ws.Range("M30").Value = 2835
$ws.Range("M31").Value = 2835
$ws.Range("M32").Value = 2835
$ws.Range("M33").Value = 2835
$ws.Range("M34").Value = 2835
$ws.Range("M35").Value = 2835
$ws.Range("M36").Value = 2835
$ws.Range("M41").Value = 2835
$ws.Range("M42").Value = 2835
$ws.Range("M43").Value = 2835
$ws.Range("M44").Value = 2835
$ws.Range("M45").Value = 2835
$ws.Range("M46").Value = 2835
$ws.Range("M47").Value = 1935
$ws.Range("M48").Value = 1935
$ws.Range("M49").Value = 1440
$ws.Range("M50").Value = 1440
$ws.Range("M51").Value = 1440
$ws.Range("M52").Value = 1440
$ws.Range("M53").Value = 1440
$ws.Range("M54").Value = 1440
$ws.Range("M55").Value = 1440
$ws.Range("M56").Value = 1440
$ws.Range("M57").Value = 1440
$ws.Range("M58").Value = 1440
$ws.Range("M59").Value = 1440
$ws.Range("M60").Value = 1440
$ws.Range("M61").Value = 1440
$ws.Range("M62").Value = 1440
$ws.Range("M63").Value = 1440
$ws.Range("M64").Value = 1440
$ws.Range("M65").Value = 1440
$ws.Range("M77").Value = 900
$ws.Range("M78").Value = 2835
$ws.Range("M79").Value = 1935
$ws.Range("M80").Value = 1440
$ws.Range("M98").Value = 900
$ws.Range("M99").Value = 900
$ws.Range("M100").Value = 900
$ws.Range("M101").Value = 900
$ws.Range("M102").Value = 900
$ws.Range("M103").Value = 900
$ws.Range("M104").Value = 900
$ws.Range("M105").Value = 900
$ws.Range("M112").Value = 1935
$ws.Range("M113").Value = 1935
$ws.Range("M114").Value = 1935
$ws.Range("M115").Value = 1935
$ws.Range("M116").Value = 1440
$ws.Range("M117").Value = 1440
$ws.Range("M118").Value = 1440
$ws.Range("M119").Value = 1440
$ws.Range("M120").Value = 1935
$ws.Range("M121").Value = 2835
$ws.Range("M122").Value = 2835
$ws.Range("M123").Value = 2835

$ws = $wb.Worksheets.Item("Wall")
$ws.Range("M2").Value = 2835

$ws = $wb.Worksheets.Item("FlowTerminal")
$ws.Range("M3").Value = 2835
$ws.Range("M4").Value = 1440

$ws = $wb.Worksheets.Item("FurnishingElement")
$ws.Range("M2").Value = 2835
